# Apply updated dSF (column F) values as described in the commit:
# "repull data, push all data, mean calculation"
# Column F ("dSF") values were re-pulled from source and differ from the
# previously copied dS0 (column E) values for a subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    5  = -6
    8  = -3
    19 = 8
    20 = 5
    21 = 2
    25 = 5
    26 = 3
    28 = 7
    29 = 6
    30 = 1
    31 = -1
    35 = -1
    36 = 1
    37 = 2
    39 = 0
    44 = 2
    46 = 1
    48 = -3
    49 = 1
    53 = 0
    54 = -1
    58 = 3
    61 = -2
    63 = 3
    64 = -9
    65 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
